$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (including the date style on column A) from the last
# existing data row (233) down to the 5 new rows being appended (234-238).
$ws.Range("A233:D233").Copy()
$ws.Range("A234:D238").PasteSpecial(-4122)

$data = @(
    @(44308, 5, 53, 309.3084330318062),
    @(44309, 3, 47, 274.2923840093376),
    @(44310, 9, 49, 285.9644003501605),
    @(44311, 2, 41, 239.276334986869),
    @(44312, 5, 35, 204.2602859644004)
)

$r = 234
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
